$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: write cell values for the new rows (19-24) describing the
# self-play / shuffle-position evaluation runs added to the log.
$ws.Range("A20").Value = "SD-SP"
$ws.Range("B20").Value = "PPO use step distance reward + multiply critic lr + train every episode + self play vs. Random"
$ws.Range("C20").Value = "能过1 2 8"
$ws.Range("D20").Value = "有时候折返能过1 2 8"
$ws.Range("E20").Value = "赢了15把，但是不太抗干扰，而且方式比较诡异"
$ws.Range("J20").Value = "贴墙走"
$ws.Range("K20").Value = "在墙上弹几下"
$ws.Range("O20").Value = "map*_use_step_dist_self_play"

$ws.Range("A21").Value = "ED-SP"
$ws.Range("B21").Value = "PPO use episode distance reward + multiply critic lr + train every episode + self play vs. Random"
$ws.Range("C21").Value = "过了好几个，1特别稳"
$ws.Range("J21").Value = "很快"
$ws.Range("K21").Value = "很快"
$ws.Range("O21").Value = "map*_use_dist_self_play"

$ws.Range("A22").Value = "SD-SP-Spos"
$ws.Range("B22").Value = "PPO use step distance reward + multiply critic lr + train every episode + self play + shuffle position vs. Random"
$ws.Range("H22").Value = "到不了最上面一层"
$ws.Range("I22").Value = "进不去"
$ws.Range("O22").Value = "map*_use_step_dist_self_play_shuffle_pos"

$ws.Range("A23").Value = "ED-SP-Spos"
$ws.Range("B23").Value = "PPO use episode distance reward + multiply critic lr + train every episode + self play + shuffle position vs. Random"
$ws.Range("C23").Value = "甚至好像学会了堵路"
$ws.Range("D23").Value = "并非稳赢，有时卡墙or折返"
$ws.Range("E23").Value = "稳的一匹"
$ws.Range("F23").Value = "稳的一匹"
$ws.Range("G23").Value = "稳的一匹，是目前最稳的"
$ws.Range("H23").Value = "过不了小房间"
$ws.Range("I23").Value = "稳的一匹，是目前最稳的"
$ws.Range("J23").Value = "很稳"
$ws.Range("K23").Value = "并非稳赢，会折返"
$ws.Range("L23").Value = "好稳啊"
$ws.Range("M23").Value = "会贴墙，走的太慢了，不过还是蛮远的"
$ws.Range("N23").Value = "打Baseline或random的冲线率都是65%，跟baseline打rewar的胜率是甚至各把11给过了一次！"
$ws.Range("O23").Value = "map*_use_dist_self_play_shuffle_pos"

$ws.Range("A24").Value = "ED-SP-Spos-3F"
$ws.Range("B24").Value = "PPO use episode distance reward + multiply critic lr + train every episode + self play + shuffle position + actor 2 layers + 3 frames vs. Random"
$ws.Range("H24").Value = "R"
$ws.Range("M24").Value = "转圈且行动缓慢"
$ws.Range("N24").Value = "懒得一批"
$ws.Range("O24").Value = "map*_use_dist_self_play_shuffle_pos_[actor]2layers_[frames]3"

# Step 2: apply cell background fills by copying the format from existing
# cells that already carry the desired style (keeps reusing the workbook's
# existing style indices instead of minting new ones).
$targets_1 = @("C20", "D20", "E20", "J20", "K20", "C21", "J21", "K21", "C23", "E23", "F23", "G23", "I23", "J23", "K23", "L23")
foreach ($t in $targets_1) {
    $ws.Range("C2").Copy()
    $ws.Range($t).PasteSpecial(-4122)
}

$targets_2 = @("H22", "I22", "H23", "M23", "M24", "N24")
foreach ($t in $targets_2) {
    $ws.Range("D2").Copy()
    $ws.Range($t).PasteSpecial(-4122)
}

$targets_3 = @("D23", "N23")
foreach ($t in $targets_3) {
    $ws.Range("L6").Copy()
    $ws.Range($t).PasteSpecial(-4122)
}

$targets_6 = @("A19", "A20", "A21", "A22", "A23", "A24")
foreach ($t in $targets_6) {
    $ws.Range("A2").Copy()
    $ws.Range($t).PasteSpecial(-4122)
}

# Step 3: view state - rezoom and move the active selection as in the
# authored edit.
$ws.Select()
$excel.ActiveWindow.Zoom = 117
$ws.Range("H14").Select()
